# Apply the cryptos list refresh (prices + 1h volume %) from the
# scraped GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the cell to keep a literal text value (matches the source
    # workbook's inlineStr cells) instead of letting Excel re-interpret
    # numeric-looking strings (e.g. "36.98", "1.00") as numbers, and
    # restore the cell's original style afterwards so no stray number
    # format / quote-prefix style is left behind.
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

# --- Price (D) / Volume(1h) (E) refresh for each coin row ---
Set-TextValue $ws.Range("D2") "51.672.64"
Set-TextValue $ws.Range("E2") "  +1.08%  "
Set-TextValue $ws.Range("D3") "3.053.14"
Set-TextValue $ws.Range("E3") "  +3.09%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.11%  "
Set-TextValue $ws.Range("D5") "384.75"
Set-TextValue $ws.Range("E5") "  +0.92%  "
Set-TextValue $ws.Range("D6") "103.01"
Set-TextValue $ws.Range("E6") "  +0.77%  "
Set-TextValue $ws.Range("E7") "  -0.05%  "
Set-TextValue $ws.Range("E8") "  +0.01%  "
Set-TextValue $ws.Range("E9") "  -0.94%  "
Set-TextValue $ws.Range("D10") "36.98"
Set-TextValue $ws.Range("E10") "  +1.35%  "
Set-TextValue $ws.Range("E11") "  +0.22%  "
Set-TextValue $ws.Range("D12") "0.0867"
Set-TextValue $ws.Range("E12") "  +1.37%  "
Set-TextValue $ws.Range("D13") "3.524.03"
Set-TextValue $ws.Range("E13") "  +3.10%  "
Set-TextValue $ws.Range("D14") "18.77"
Set-TextValue $ws.Range("E14") "  +2.21%  "
Set-TextValue $ws.Range("D15") "7.75"
Set-TextValue $ws.Range("E15") "  -0.63%  "
Set-TextValue $ws.Range("D16") "3.038.28"
Set-TextValue $ws.Range("E16") "  +2.90%  "
Set-TextValue $ws.Range("D17") "0.981"
Set-TextValue $ws.Range("E17") "  -1.79%  "
Set-TextValue $ws.Range("D18") "10.56"
Set-TextValue $ws.Range("E18") "  -7.93%  "
Set-TextValue $ws.Range("D19") "51.674.30"
Set-TextValue $ws.Range("E19") "  +0.95%  "
Set-TextValue $ws.Range("E21") "  +0.61%  "
Set-TextValue $ws.Range("D22") "0.0₃0966"
Set-TextValue $ws.Range("E22") "  +0.33%  "
Set-TextValue $ws.Range("D23") "70.06"
Set-TextValue $ws.Range("E23") "  -0.23%  "
Set-TextValue $ws.Range("D24") "267.42"
Set-TextValue $ws.Range("E24") "  +0.17%  "
Set-TextValue $ws.Range("D25") "3.18"
Set-TextValue $ws.Range("E25") "  -4.38%  "
Set-TextValue $ws.Range("D26") "8.46"
Set-TextValue $ws.Range("E26") "  +7.21%  "
Set-TextValue $ws.Range("D29") "26.49"
Set-TextValue $ws.Range("E29") "  +2.38%  "
Set-TextValue $ws.Range("E30") "  +0.00%  "
Set-TextValue $ws.Range("E31") "  -3.17%  "
Set-TextValue $ws.Range("D32") "10.30"
Set-TextValue $ws.Range("E32") "  -0.12%  "
Set-TextValue $ws.Range("D33") "34.23"
Set-TextValue $ws.Range("E33") "  -0.36%  "
Set-TextValue $ws.Range("E34") "  +0.21%  "
Set-TextValue $ws.Range("D35") "50.52"
Set-TextValue $ws.Range("E35") "  -1.17%  "
Set-TextValue $ws.Range("D36") "0.0444"
Set-TextValue $ws.Range("E36") "  +1.81%  "
Set-TextValue $ws.Range("E37") "  -0.07%  "
Set-TextValue $ws.Range("E38") "  +4.01%  "
Set-TextValue $ws.Range("E39") "  +5.22%  "
Set-TextValue $ws.Range("D40") "17.15"
Set-TextValue $ws.Range("E40") "  +3.69%  "
Set-TextValue $ws.Range("E41") "  +2.36%  "
Set-TextValue $ws.Range("D42") "127.77"
Set-TextValue $ws.Range("E42") "  +2.41%  "
Set-TextValue $ws.Range("E43") "  -0.31%  "
Set-TextValue $ws.Range("E44") "  +0.75%  "
Set-TextValue $ws.Range("D45") "3.69"
Set-TextValue $ws.Range("E45") "  +4.39%  "
Set-TextValue $ws.Range("D46") "21.73"
Set-TextValue $ws.Range("E46") "  +1.33%  "
Set-TextValue $ws.Range("D47") "2.52"
Set-TextValue $ws.Range("E47") "  +6.43%  "
Set-TextValue $ws.Range("E48") "  +3.10%  "
Set-TextValue $ws.Range("D49") "2.037.58"
Set-TextValue $ws.Range("E49") "  -0.63%  "
Set-TextValue $ws.Range("D50") "3.348.05"
Set-TextValue $ws.Range("E50") "  +2.90%  "
Set-TextValue $ws.Range("E51") "  +6.84%  "

# --- Rows 27/28 swapped: Kaspa now ranks above RenderToken ---
Set-TextValue $ws.Range("B27") "Kaspa"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D27") "0.173"
Set-TextValue $ws.Range("E27") "  +4.15%  "

Set-TextValue $ws.Range("B28") "RenderToken"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D28") "7.34"
Set-TextValue $ws.Range("E28") "  +1.60%  "

